$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 89, shifting the existing rows 89-96 down to 90-97.
$ws.Rows.Item(89).Insert()

# Populate the new row 89 with the new weekly price record.
$ws.Range("A89").Value = 4
$ws.Range("B89").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C89").Value = "Los Lagos"
$ws.Range("D89").Value = 45267
$ws.Range("E89").Value = 10
$ws.Range("F89").Value = 300000000
$ws.Range("G89").Value = "Espárragos"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 100
$ws.Range("K89").Value = 2000
$ws.Range("L89").Value = 2000
$ws.Range("M89").Value = 2000
$ws.Range("N89").Value = "`$/kilo"
$ws.Range("O89").Value = "Provincia de Linares"
$ws.Range("P89").Value = 2000
$ws.Range("Q89").Value = 1
$ws.Range("R89").Value = "Hortaliza"
